$d = $word.ActiveDocument

$d.Content.Find.Execute("5+92=", $true, $false, $false, $false, $false, $true, 1, $false, "25+45=", 2) | Out-Null
$d.Content.Find.Execute("44+42=", $true, $false, $false, $false, $false, $true, 1, $false, "26+40=", 2) | Out-Null
$d.Content.Find.Execute("27+68=", $true, $false, $false, $false, $false, $true, 1, $false, "46+45=", 2) | Out-Null
$d.Content.Find.Execute("72-4=", $true, $false, $false, $false, $false, $true, 1, $false, "24+70=", 2) | Out-Null
$d.Content.Find.Execute("66+32=", $true, $false, $false, $false, $false, $true, 1, $false, "5+37=", 2) | Out-Null
$d.Content.Find.Execute("4+6=", $true, $false, $false, $false, $false, $true, 1, $false, "42+21=", 2) | Out-Null
$d.Content.Find.Execute("76-30=", $true, $false, $false, $false, $false, $true, 1, $false, "97-59=", 2) | Out-Null
$d.Content.Find.Execute("99-30=", $true, $false, $false, $false, $false, $true, 1, $false, "8+51=", 2) | Out-Null
$d.Content.Find.Execute("72+18=", $true, $false, $false, $false, $false, $true, 1, $false, "38+1=", 2) | Out-Null
$d.Content.Find.Execute("54+24=", $true, $false, $false, $false, $false, $true, 1, $false, "21+52=", 2) | Out-Null
$d.Content.Find.Execute("57-18=", $true, $false, $false, $false, $false, $true, 1, $false, "46-29=", 2) | Out-Null
$d.Content.Find.Execute("57-34=", $true, $false, $false, $false, $false, $true, 1, $false, "30-12=", 2) | Out-Null
$d.Content.Find.Execute("93-64=", $true, $false, $false, $false, $false, $true, 1, $false, "41+15=", 2) | Out-Null
$d.Content.Find.Execute("43-19=", $true, $false, $false, $false, $false, $true, 1, $false, "6+22=", 2) | Out-Null
$d.Content.Find.Execute("21+14=", $true, $false, $false, $false, $false, $true, 1, $false, "12+73=", 2) | Out-Null
$d.Content.Find.Execute("84-30=", $true, $false, $false, $false, $false, $true, 1, $false, "19+79=", 2) | Out-Null
$d.Content.Find.Execute("68-20=", $true, $false, $false, $false, $false, $true, 1, $false, "16-15=", 2) | Out-Null
$d.Content.Find.Execute("62-62=", $true, $false, $false, $false, $false, $true, 1, $false, "21-9=", 2) | Out-Null
$d.Content.Find.Execute("65-64=", $true, $false, $false, $false, $false, $true, 1, $false, "44+13=", 2) | Out-Null
$d.Content.Find.Execute("42-34=", $true, $false, $false, $false, $false, $true, 1, $false, "46+37=", 2) | Out-Null
$d.Content.Find.Execute("89-53=", $true, $false, $false, $false, $false, $true, 1, $false, "32-21=", 2) | Out-Null
$d.Content.Find.Execute("83-59=", $true, $false, $false, $false, $false, $true, 1, $false, "42+55=", 2) | Out-Null
$d.Content.Find.Execute("12+76=", $true, $false, $false, $false, $false, $true, 1, $false, "6+17=", 2) | Out-Null
$d.Content.Find.Execute("90-84=", $true, $false, $false, $false, $false, $true, 1, $false, "11-7=", 2) | Out-Null
$d.Content.Find.Execute("87+1=", $true, $false, $false, $false, $false, $true, 1, $false, "71-20=", 2) | Out-Null
$d.Content.Find.Execute("79-15=", $true, $false, $false, $false, $false, $true, 1, $false, "9+84=", 2) | Out-Null
$d.Content.Find.Execute("4+90=", $true, $false, $false, $false, $false, $true, 1, $false, "32-29=", 2) | Out-Null
$d.Content.Find.Execute("2+42=", $true, $false, $false, $false, $false, $true, 1, $false, "33-20=", 2) | Out-Null
$d.Content.Find.Execute("51+25=", $true, $false, $false, $false, $false, $true, 1, $false, "98-26=", 2) | Out-Null
$d.Content.Find.Execute("3+69=", $true, $false, $false, $false, $false, $true, 1, $false, "27+33=", 2) | Out-Null
$d.Content.Find.Execute("85-33=", $true, $false, $false, $false, $false, $true, 1, $false, "69+10=", 2) | Out-Null
$d.Content.Find.Execute("13+23=", $true, $false, $false, $false, $false, $true, 1, $false, "96-75=", 2) | Out-Null
$d.Content.Find.Execute("84-59=", $true, $false, $false, $false, $false, $true, 1, $false, "78+4=", 2) | Out-Null
$d.Content.Find.Execute("11-6=", $true, $false, $false, $false, $false, $true, 1, $false, "97-36=", 2) | Out-Null
$d.Content.Find.Execute("20+66=", $true, $false, $false, $false, $false, $true, 1, $false, "69+24=", 2) | Out-Null
$d.Content.Find.Execute("99-8=", $true, $false, $false, $false, $false, $true, 1, $false, "11+20=", 2) | Out-Null
$d.Content.Find.Execute("92-86=", $true, $false, $false, $false, $false, $true, 1, $false, "64-13=", 2) | Out-Null
$d.Content.Find.Execute("93-75=", $true, $false, $false, $false, $false, $true, 1, $false, "48+51=", 2) | Out-Null
$d.Content.Find.Execute("61-28=", $true, $false, $false, $false, $false, $true, 1, $false, "79-48=", 2) | Out-Null
$d.Content.Find.Execute("72-59=", $true, $false, $false, $false, $false, $true, 1, $false, "96-6=", 2) | Out-Null
$d.Content.Find.Execute("59-46=", $true, $false, $false, $false, $false, $true, 1, $false, "44+20=", 2) | Out-Null
$d.Content.Find.Execute("68-53=", $true, $false, $false, $false, $false, $true, 1, $false, "30+36=", 2) | Out-Null
$d.Content.Find.Execute("10+28=", $true, $false, $false, $false, $false, $true, 1, $false, "90-76=", 2) | Out-Null
$d.Content.Find.Execute("90-38=", $true, $false, $false, $false, $false, $true, 1, $false, "42-8=", 2) | Out-Null
$d.Content.Find.Execute("37-34=", $true, $false, $false, $false, $false, $true, 1, $false, "19+18=", 2) | Out-Null
$d.Content.Find.Execute("50-16=", $true, $false, $false, $false, $false, $true, 1, $false, "36+9=", 2) | Out-Null
$d.Content.Find.Execute("39+18=", $true, $false, $false, $false, $false, $true, 1, $false, "57+41=", 2) | Out-Null
$d.Content.Find.Execute("13+9=", $true, $false, $false, $false, $false, $true, 1, $false, "28-28=", 2) | Out-Null
$d.Content.Find.Execute("50+47=", $true, $false, $false, $false, $false, $true, 1, $false, "22+7=", 2) | Out-Null
$d.Content.Find.Execute("80+3=", $true, $false, $false, $false, $false, $true, 1, $false, "73-6=", 2) | Out-Null
$d.Content.Find.Execute("47-20=", $true, $false, $false, $false, $false, $true, 1, $false, "30+57=", 2) | Out-Null
$d.Content.Find.Execute("10+64=", $true, $false, $false, $false, $false, $true, 1, $false, "64+35=", 2) | Out-Null
$d.Content.Find.Execute("1+97=", $true, $false, $false, $false, $false, $true, 1, $false, "55-43=", 2) | Out-Null
$d.Content.Find.Execute("16+32=", $true, $false, $false, $false, $false, $true, 1, $false, "67-2=", 2) | Out-Null
$d.Content.Find.Execute("77-35=", $true, $false, $false, $false, $false, $true, 1, $false, "87-72=", 2) | Out-Null
$d.Content.Find.Execute("33-4=", $true, $false, $false, $false, $false, $true, 1, $false, "74+12=", 2) | Out-Null
$d.Content.Find.Execute("80+10=", $true, $false, $false, $false, $false, $true, 1, $false, "35+40=", 2) | Out-Null
$d.Content.Find.Execute("51+44=", $true, $false, $false, $false, $false, $true, 1, $false, "68-58=", 2) | Out-Null
$d.Content.Find.Execute("73+7=", $true, $false, $false, $false, $false, $true, 1, $false, "93-78=", 2) | Out-Null
$d.Content.Find.Execute("5+76=", $true, $false, $false, $false, $false, $true, 1, $false, "40+22=", 2) | Out-Null
$d.Content.Find.Execute("30+61=", $true, $false, $false, $false, $false, $true, 1, $false, "3+34=", 2) | Out-Null
$d.Content.Find.Execute("13+68=", $true, $false, $false, $false, $false, $true, 1, $false, "7+41=", 2) | Out-Null
$d.Content.Find.Execute("70+4=", $true, $false, $false, $false, $false, $true, 1, $false, "40-39=", 2) | Out-Null
$d.Content.Find.Execute("77-57=", $true, $false, $false, $false, $false, $true, 1, $false, "16-1=", 2) | Out-Null
$d.Content.Find.Execute("49+41=", $true, $false, $false, $false, $false, $true, 1, $false, "71-49=", 2) | Out-Null
$d.Content.Find.Execute("35+10=", $true, $false, $false, $false, $false, $true, 1, $false, "97-31=", 2) | Out-Null
$d.Content.Find.Execute("89-34=", $true, $false, $false, $false, $false, $true, 1, $false, "42+41=", 2) | Out-Null
$d.Content.Find.Execute("28+13=", $true, $false, $false, $false, $false, $true, 1, $false, "9+43=", 2) | Out-Null
$d.Content.Find.Execute("61-51=", $true, $false, $false, $false, $false, $true, 1, $false, "37+21=", 2) | Out-Null
$d.Content.Find.Execute("14-4=", $true, $false, $false, $false, $false, $true, 1, $false, "71-52=", 2) | Out-Null
$d.Content.Find.Execute("4+28=", $true, $false, $false, $false, $false, $true, 1, $false, "53-14=", 2) | Out-Null
$d.Content.Find.Execute("40-36=", $true, $false, $false, $false, $false, $true, 1, $false, "7+4=", 2) | Out-Null
$d.Content.Find.Execute("24+53=", $true, $false, $false, $false, $false, $true, 1, $false, "72-15=", 2) | Out-Null
$d.Content.Find.Execute("50-0=", $true, $false, $false, $false, $false, $true, 1, $false, "60-24=", 2) | Out-Null
$d.Content.Find.Execute("49-12=", $true, $false, $false, $false, $false, $true, 1, $false, "8+43=", 2) | Out-Null
$d.Content.Find.Execute("50+29=", $true, $false, $false, $false, $false, $true, 1, $false, "44+9=", 2) | Out-Null
$d.Content.Find.Execute("79-40=", $true, $false, $false, $false, $false, $true, 1, $false, "25-24=", 2) | Out-Null
$d.Content.Find.Execute("49+2=", $true, $false, $false, $false, $false, $true, 1, $false, "7+20=", 2) | Out-Null
$d.Content.Find.Execute("83-37=", $true, $false, $false, $false, $false, $true, 1, $false, "77-72=", 2) | Out-Null
$d.Content.Find.Execute("28+9=", $true, $false, $false, $false, $false, $true, 1, $false, "44+23=", 2) | Out-Null
$d.Content.Find.Execute("17+53=", $true, $false, $false, $false, $false, $true, 1, $false, "38+56=", 2) | Out-Null
$d.Content.Find.Execute("56+28=", $true, $false, $false, $false, $false, $true, 1, $false, "34+32=", 2) | Out-Null
$d.Content.Find.Execute("27+71=", $true, $false, $false, $false, $false, $true, 1, $false, "19+27=", 2) | Out-Null
$d.Content.Find.Execute("75-73=", $true, $false, $false, $false, $false, $true, 1, $false, "85-34=", 2) | Out-Null
$d.Content.Find.Execute("55+40=", $true, $false, $false, $false, $false, $true, 1, $false, "8+89=", 2) | Out-Null
$d.Content.Find.Execute("63-38=", $true, $false, $false, $false, $false, $true, 1, $false, "13+55=", 2) | Out-Null
$d.Content.Find.Execute("37+38=", $true, $false, $false, $false, $false, $true, 1, $false, "54-47=", 2) | Out-Null
$d.Content.Find.Execute("45-20=", $true, $false, $false, $false, $false, $true, 1, $false, "41-17=", 2) | Out-Null
$d.Content.Find.Execute("27+61=", $true, $false, $false, $false, $false, $true, 1, $false, "30+30=", 2) | Out-Null
$d.Content.Find.Execute("83-2=", $true, $false, $false, $false, $false, $true, 1, $false, "61+10=", 2) | Out-Null
$d.Content.Find.Execute("89-26=", $true, $false, $false, $false, $false, $true, 1, $false, "95-38=", 2) | Out-Null
$d.Content.Find.Execute("55-0=", $true, $false, $false, $false, $false, $true, 1, $false, "56-2=", 2) | Out-Null
$d.Content.Find.Execute("66-51=", $true, $false, $false, $false, $false, $true, 1, $false, "92-31=", 2) | Out-Null
$d.Content.Find.Execute("15+30=", $true, $false, $false, $false, $false, $true, 1, $false, "51+5=", 2) | Out-Null
$d.Content.Find.Execute("9+14=", $true, $false, $false, $false, $false, $true, 1, $false, "5+40=", 2) | Out-Null
$d.Content.Find.Execute("99-65=", $true, $false, $false, $false, $false, $true, 1, $false, "95-69=", 2) | Out-Null
$d.Content.Find.Execute("40+29=", $true, $false, $false, $false, $false, $true, 1, $false, "48+0=", 2) | Out-Null
$d.Content.Find.Execute("36+59=", $true, $false, $false, $false, $false, $true, 1, $false, "82-51=", 2) | Out-Null
$d.Content.Find.Execute("53-34=", $true, $false, $false, $false, $false, $true, 1, $false, "27+39=", 2) | Out-Null
$d.Content.Find.Execute("87-54=", $true, $false, $false, $false, $false, $true, 1, $false, "96-87=", 2) | Out-Null

Write-Host "Done replacing 100 expressions"
